$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.32
$ws.Range("G2").Value = 1.39
$ws.Range("H2").Value = 7
$ws.Range("G3").Value = 9.199999999999999
$ws.Range("I3").Value = 1.57
$ws.Range("K3").Value = 4.7
$ws.Range("N3").Value = 3.5
$ws.Range("P3").Value = 1.86
$ws.Range("Q3").Value = 2.06
$ws.Range("R3").Value = 1.31
$ws.Range("S3").Value = 3.85
$ws.Range("T3").Value = 2.2
$ws.Range("U3").Value = 1.73
$ws.Range("Y3").Value = 7.4
$ws.Range("AC3").Value = 10.5
$ws.Range("AE3").Value = 18.5
$ws.Range("AK3").Value = 1000
$ws.Range("I4").Value = 2.88
$ws.Range("P4").Value = 1.69
$ws.Range("F5").Value = 2.58
$ws.Range("G5").Value = 3.05
$ws.Range("H5").Value = 3.05
$ws.Range("I5").Value = 3.7
$ws.Range("J5").Value = 2.94
$ws.Range("K5").Value = 3.3
$ws.Range("P5").Value = 1.56
$ws.Range("Q5").Value = 2.3
$ws.Range("F6").Value = 2.24
$ws.Range("G6").Value = 2.48
$ws.Range("I6").Value = 4.3
$ws.Range("F7").Value = 4
$ws.Range("I7").Value = 2.18
$ws.Range("K7").Value = 4
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 2.3
$ws.Range("H8").Value = 3.9
$ws.Range("I8").Value = 4.7
$ws.Range("J8").Value = 3.25
$ws.Range("K8").Value = 4.2
$ws.Range("P8").Value = 1.73
$ws.Range("Q8").Value = 2.12
$ws.Range("F9").Value = 2.04
$ws.Range("G9").Value = 2.54
$ws.Range("H9").Value = 3.35
$ws.Range("I9").Value = 4.7
$ws.Range("J9").Value = 3.25
$ws.Range("K9").Value = 4.8
$ws.Range("P9").Value = 1.61
$ws.Range("Q9").Value = 2.06
$ws.Range("F10").Value = 1.55
$ws.Range("H10").Value = 5.4
$ws.Range("J10").Value = 4.6
$ws.Range("G11").Value = 2.46
$ws.Range("I11").Value = 3.75
$ws.Range("J11").Value = 3.5
$ws.Range("K11").Value = 4.2
$ws.Range("Q11").Value = 1.72
$ws.Range("F12").Value = 2.72
$ws.Range("G12").Value = 2.76
$ws.Range("I12").Value = 2.84
$ws.Range("K12").Value = 3.55
$ws.Range("Q12").Value = 1.91
$ws.Range("R12").Value = 1.4
$ws.Range("S12").Value = 3.3
$ws.Range("T12").Value = 1.72
$ws.Range("Y12").Value = 14
$ws.Range("AC12").Value = 8.199999999999999
$ws.Range("AI12").Value = 44
$ws.Range("AJ12").Value = 48
$ws.Range("AM12").Value = 100
$ws.Range("AN12").Value = 32
$ws.Range("G13").Value = 1.92
$ws.Range("H13").Value = 3.95
$ws.Range("I13").Value = 6.4
$ws.Range("J13").Value = 3.05
$ws.Range("P13").Value = 1.53
$ws.Range("P14").Value = 1.25
$ws.Range("Q14").Value = 1.01
$ws.Range("G15").Value = 2.34
$ws.Range("H15").Value = 3.1
